$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 11 (the "P11 total" row), pushing the
# total row and the footer row down by one. This makes room for a new item
# (item #4 "SUGARLO PLUS..."), while the previous item #4
# ("صابون ديتول اوريجنيال 115 جم") becomes item #5 on the new row 11.
$ws.Rows.Item(11).Insert()
$ws.Rows.Item(11).RowHeight = 25.5

# Copy the formatting of row 10 (the item row just above) onto the freshly
# inserted row 11 so it matches the other item rows exactly.
$ws.Range("A10:Q10").Copy()
$ws.Range("A11:Q11").PasteSpecial(-4122)

# Re-create the merges for the new row 11, matching the pattern used by the
# other item rows (7-10).
$ws.Range("A11:B11").Merge()
$ws.Range("C11:G11").Merge()
$ws.Range("H11:K11").Merge()
$ws.Range("L11:M11").Merge()
$ws.Range("N11:O11").Merge()

# Row 11 now holds what used to be item #4: bump its item number to 5 and
# copy over the old values (name/balance/limit/price/sale price/count).
$ws.Range("A11").Value = 5
$ws.Range("C11").Value = "صابون ديتول اوريجنيال 115 جم"
$ws.Range("H11").Value = "0:0"
$ws.Range("L11").Value = "0"
$ws.Range("N11").Value = "30.00"
$ws.Range("P11").Value = "30.0000"
$ws.Range("Q11").Value = "1:0"

# Row 10 becomes the new item #4: "SUGARLO PLUS 50/1000MG 30 F.C. TABS".
$ws.Range("C10").Value = "SUGARLO PLUS 50/1000MG 30 F.C. TABS"
$ws.Range("H10").Value = "1:0"
$ws.Range("L10").Value = "1"
$ws.Range("N10").Value = "136.50"
$ws.Range("P10").Value = "45.0450"
$ws.Range("Q10").Value = "0:1"

# Update the running total (old P11, now P12) to include the new item's sale
# price.
$ws.Range("P12").Value = 186.04

# Update the generated timestamp in the footer (was A12, now shifted to A13).
$ws.Range("A13").Value = "Monday, 28 July, 2025 10:09 AM"
